$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" is column G. Determine the last used row so we scan the full table.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -ne $value -and $value -is [string] -and $value.Contains(",")) {
        # Split the comma separated list of recorders, e.g. "dnasr281@gmail.com, System"
        $parts = $value -split ",\s*"

        # Find the entry that is exactly "System" (case sensitive - keep lowercase "system" untouched)
        $idx = -1
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($parts[$i] -ceq "System") {
                $idx = $i
                break
            }
        }

        if ($idx -ge 0) {
            # Move "System" to the front, preserving the relative order of the remaining entries
            $newParts = @("System")
            for ($i = 0; $i -lt $parts.Length; $i++) {
                if ($i -ne $idx) {
                    $newParts += $parts[$i]
                }
            }
            $newValue = [string]::Join(", ", $newParts)

            if ($newValue -ne $value) {
                $cell.Value2 = $newValue
            }
        }
    }
}
